$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old content in rows 2-38 across columns B-F (rows 36-38 will be removed entirely)
$ws.Range("B2:F38").ClearContents()

# Delete rows 36-38 (shrinks used range / dimension from F38 to F35)
$ws.Range("A36:F38").Delete()

# Set new ticker values for columns B and C, rows 2-35
$ws.Range("B2").Value = "NSE:ASIANENE"
$ws.Range("C2").Value = "NSE:AAVAS"
$ws.Range("B3").Value = "NSE:CARTRADE"
$ws.Range("C3").Value = "NSE:AKASH"
$ws.Range("B4").Value = "NSE:CHOICEIN"
$ws.Range("C4").Value = "NSE:AWL"
$ws.Range("B5").Value = "NSE:IVC"
$ws.Range("C5").Value = "NSE:BBOX"
$ws.Range("B6").Value = "NSE:JINDALSAW"
$ws.Range("C6").Value = "NSE:BHAGCHEM"
$ws.Range("B7").Value = "NSE:LUPIN"
$ws.Range("C7").Value = "NSE:BHARATWIRE"
$ws.Range("B8").Value = "NSE:NETWEB"
$ws.Range("C8").Value = "NSE:BIKAJI"
$ws.Range("B9").Value = "NSE:POLYCAB"
$ws.Range("C9").Value = "NSE:BORORENEW"
$ws.Range("C10").Value = "NSE:BSOFT"
$ws.Range("C11").Value = "NSE:COMPUSOFT"
$ws.Range("C12").Value = "NSE:DCAL"
$ws.Range("C13").Value = "NSE:DEN"
$ws.Range("C14").Value = "NSE:DENORA"
$ws.Range("C15").Value = "NSE:EMMBI"
$ws.Range("C16").Value = "NSE:FOSECOIND"
$ws.Range("C17").Value = "NSE:HARRMALAYA"
$ws.Range("C18").Value = "NSE:ISGEC"
$ws.Range("C19").Value = "NSE:ISMTLTD"
$ws.Range("C20").Value = "NSE:JAYBARMARU"
$ws.Range("C21").Value = "NSE:MANOMAY"
$ws.Range("C22").Value = "NSE:MANYAVAR"
$ws.Range("C23").Value = "NSE:MITCON"
$ws.Range("C24").Value = "NSE:MTNL"
$ws.Range("C25").Value = "NSE:MUNJALSHOW"
$ws.Range("C26").Value = "NSE:MURUDCERA"
$ws.Range("C27").Value = "NSE:NAGAFERT"
$ws.Range("C28").Value = "NSE:NAGREEKEXP"
$ws.Range("C29").Value = "NSE:NESCO"
$ws.Range("C30").Value = "NSE:NGIL"
$ws.Range("C31").Value = "NSE:NIPPOBATRY"
$ws.Range("C32").Value = "NSE:OBCL"
$ws.Range("C33").Value = "NSE:PATANJALI"
$ws.Range("C34").Value = "NSE:RPGLIFE"
$ws.Range("C35").Value = "NSE:RSWM"
